$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8 from 45233 to 45243
$ws.Range("C2:C8").Value = 45243
